# Commit: "replace GR with Gr in fix.php"
# The underlying XML diff shows every shared string that contains the
# substring "GR" (always followed by a group number, e.g. "GR1", "GR2", ...)
# had it changed to "Gr" (e.g. "DS-GR1 Abdul Aziz" -> "DS-Gr1 Abdul Aziz").
# Strings that already used "Gr" (e.g. "OR Gr1 Aasma") are untouched.
#
# Apply this as a case-sensitive whole-workbook find & replace of "GR" -> "Gr"
# across every worksheet, so every cell (the schedule appears on all 5 day
# sheets, reusing the same course/group labels) gets updated consistently.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Params: What, Replacement, LookAt:=xlPart(2), SearchOrder:=xlByRows(1),
    # MatchCase:=$true, MatchByte:=$false, SearchFormat:=$false
    $ws.Cells.Replace("GR", "Gr", 2, 1, $true, $false, $false)
}
